# PlayerPerformance_5927.xlsx edit script
# ------------------------------------------------------------------
# Goal (per commit message "added code to scrape more data about a
# player's batting performance in a match, also updated the excel
# sheets"):
#   1. Insert a new first worksheet "Player Info" holding the
#      player's biographical info (ID/NAME/BATTING_HAND/BOWL_STYLE).
#   2. Keep the existing batting log as the second worksheet renamed
#      "ODI Batting", but:
#        - rename column D header MATCH_CARD_LINK -> MATCH_CODE
#        - replace the full scorecard URL in column D with just the
#          numeric match code
#        - append a newly-scraped match row (row 11)
# ------------------------------------------------------------------

# Every column in this workbook stores its data as text, even when
# the text happens to look like a number (e.g. match codes, "10" for
# an inning number, ...). Plain `$range.Value = "123"` auto-coerces
# to a numeric cell here, so route such assignments through a text
# number-format long enough to pin the type, then drop back to the
# default ("Normal") style so no stray formatting is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Step 1: duplicate the current (only) sheet so we end up with
#     two sheets; the copy -- placed right after the original --
#     becomes the new "ODI Batting" sheet, inheriting the original
#     header styling/content verbatim (so we only have to touch the
#     handful of cells that actually changed).
$ws1.Copy($null, $ws1)
$wsBatting = $wb.Worksheets.Item(2)
$wsBatting.Name = "ODI Batting New"

# --- Step 2: update the (new) "ODI Batting" sheet -----------------
$wsBatting.Range("D1").Value = "MATCH_CODE"

Set-TextValue $wsBatting.Range("D2") "4452"
Set-TextValue $wsBatting.Range("D3") "4453"
Set-TextValue $wsBatting.Range("D4") "4563"
Set-TextValue $wsBatting.Range("D5") "4566"
Set-TextValue $wsBatting.Range("D6") "4568"
Set-TextValue $wsBatting.Range("D7") "4605"
Set-TextValue $wsBatting.Range("D8") "4608"
Set-TextValue $wsBatting.Range("D9") "4614"
Set-TextValue $wsBatting.Range("D10") "4735"

# newly scraped match row
Set-TextValue $wsBatting.Range("A11") "10"
Set-TextValue $wsBatting.Range("B11") "10"
$wsBatting.Range("C11").Value = "31/03/2023"
Set-TextValue $wsBatting.Range("D11") "4745"
$wsBatting.Range("E11").Value = "2nd"
$wsBatting.Range("F11").Value = "Sri Lanka"
$wsBatting.Range("G11").Value = "Seddon Park"
$wsBatting.Range("H11").Value = "not out"
$wsBatting.Range("I11").Value = "86*"
Set-TextValue $wsBatting.Range("J11") "113"

# --- Step 3: turn the original sheet into "Player Info" -----------
# Clear out the old batting columns E:J (header row only has 4
# columns now) and all the old data rows, keeping A1:D1's existing
# bold/bordered header style in place for reuse.
$ws1.Range("E1:J1").Clear()
$ws1.Range("A2:J10").Clear()

$ws1.Range("A1").Value = "ID"
$ws1.Range("B1").Value = "NAME"
$ws1.Range("C1").Value = "BATTING_HAND"
$ws1.Range("D1").Value = "BOWL_STYLE"

Set-TextValue $ws1.Range("A2") "5927"
$ws1.Range("B2").Value = "William Alexander Young"
$ws1.Range("C2").Value = "Right Handed"
$ws1.Range("D2").Value = "Right Arm Off Break"

$ws1.Name = "Player Info"
$wsBatting.Name = "ODI Batting"

# Restore the original "first sheet is active" selection state (the
# sheet Copy() earlier switched focus to the new second sheet).
$ws1.Activate()

Write-Output "Player Info + ODI Batting sheets updated"
